# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet.
# All data rows (2-26) currently show "2025-11-21 01:19:24" in column A and
# need to be refreshed to "2025-11-21 01:49:03" as plain text, matching the
# original string formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-21 01:19:24"
$newValue = "2025-11-21 01:49:03"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
